# Update the "Current balance" ratios for three products and refresh the
# report-generated timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BETADERM 0.1% CREAM 15 GM  (row 8)
$ws.Range("H8").Value = "6:0"

# KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF.  (row 13)
$ws.Range("H13").Value = "3:1"

# VOLTAREN 75MG/3ML 3 AMP.  (row 16)
$ws.Range("H16").Value = "6:3"

# Footer timestamp ("Friday, 1 August, 2025 4:04 PM" -> "... 4:15 PM")
$ws.Range("A23").Value = "Friday, 1 August, 2025 4:15 PM"
